# Adds three new worksheets to the workbook, matching the target OOXML:
#   - CypherOutput_Message : copy of the existing 'Message' log block
#   - StatOutput           : new, empty worksheet (placeholder for future output)
#   - StatOutput_Message   : log block for the Stat query run, followed by a second
#                            block documenting the 'empty Cypher query' validation error

$wb = $excel.ActiveWorkbook

# The ten-line message block shared by the Message / CypherOutput_Message sheets
# (connection info + the Cypher query text + the output file path).
$messageLines = @(
    'Neo4j_URL:',
    'bolt://ncias-q2251-c.nci.nih.gov:7687',
    'User_name:',
    'neo4j',
    'PWD:',
    'icdcDBneo4j0',
    'Cypher:',
    'MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE s.clinical_study_designation IN [''COTC007B'',''NCATS-COP01''] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'''') AS `Case ID` , coalesce(s.clinical_study_designation,'''') AS `Study Code` , coalesce(s.clinical_study_type,'''') AS  `Study Type`, coalesce(demo.breed,'''') AS Breed , coalesce(diag.disease_term,'''') AS Diagnosis , coalesce(diag.stage_of_disease,'''') AS `Stage of Disease` ,  coalesce(demo.patient_age_at_enrollment,'''') AS Age , coalesce(demo.sex,'''') AS Sex , coalesce(demo.neutered_indicator,'''') AS  `Neutered Status`',
    'Output:',
    'C:\Users\radhakrishnang2\Desktop\Commons_Automation\OutputFiles\TC03_Canine_Filter_Study-ALL_Neo4jData.xlsx'
)

$cypherEmptyErrorLine = 'Cypher query should not be an empty string'

# --- CypherOutput_Message: identical content to the 'Message' sheet ---
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$cypherOutputMessage = $wb.Worksheets.Add($null, $lastSheet)
$cypherOutputMessage.Name = "CypherOutput_Message"
for ($i = 0; $i -lt $messageLines.Count; $i++) {
    $cypherOutputMessage.Cells.Item($i + 1, 1).Value = $messageLines[$i]
}

# --- StatOutput: new, empty worksheet ---
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$statOutput = $wb.Worksheets.Add($null, $lastSheet)
$statOutput.Name = "StatOutput"

# --- StatOutput_Message: the standard message block, followed by the
#     'empty Cypher query' validation-error message, followed by a second
#     message block where the Cypher query text is blank ---
$count = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($count)
$statOutputMessage = $wb.Worksheets.Add($null, $lastSheet)
$statOutputMessage.Name = "StatOutput_Message"
for ($i = 0; $i -lt $messageLines.Count; $i++) {
    $statOutputMessage.Cells.Item($i + 1, 1).Value = $messageLines[$i]
}
$statOutputMessage.Cells.Item(11, 1).Value = $cypherEmptyErrorLine
for ($i = 0; $i -lt 7; $i++) {
    $statOutputMessage.Cells.Item(12 + $i, 1).Value = $messageLines[$i]
}
# Row 19 holds the (blank) Cypher query text for this second, failing run. A plain
# Value = '' clears the cell instead of leaving a real empty-text cell behind, so use
# the classic "lone apostrophe" quote-prefix trick -- Excel stores/reads that back as
# an empty string cell (still text-typed) rather than removing the cell entirely.
$statOutputMessage.Cells.Item(19, 1).Value = "'"
$statOutputMessage.Cells.Item(20, 1).Value = $messageLines[8]
$statOutputMessage.Cells.Item(21, 1).Value = $messageLines[9]

# Adding sheets shifts Excel's active tab to the last one created; restore the
# original active sheet (CypherOutput, the first tab) to match the workbook's
# prior selection state.
$wb.Worksheets.Item(1).Activate()

